$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Clone the formatting ("Ruim" style used by every data row) from row 32 down
# onto the four new rows (33-36) plus the trailing blank row (37), then fill
# in the new skill data. Copy/PasteSpecial(formats) reuses the existing
# cellXf (index 3) instead of minting a new one.
$ws.Range("A32:E32").Copy()
$ws.Range("A33:E36").PasteSpecial(-4122)
$ws.Range("A32").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "SwordAtack"
$ws.Range("C33").Value = "BasicAtack"
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0

$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "EquipSwordAtack"
$ws.Range("C34").Value = "PassiveSkill"
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0

$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "BluntStrike"
$ws.Range("C35").Value = "BasicAtack"
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0

$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "BluntAttack"
$ws.Range("C36").Value = "PassiveSkill"
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0

# Row 37 is left blank (just the carried-over formatting on A37).

$ws.Range("F36").Select()
